# Eliminate RTL Freeze terminology in spreadsheet
#
# Updates copyright/license notice, replaces "Functional RTL Freeze"
# terminology with "TRL-5" wording, rewords a couple of "RTL freeze" /
# "RTL Freeze" mentions to "release", and removes the personal e-mail
# address (with its mailto: hyperlink) that was recorded as the
# "Owner" of the "Version clearly identified" sign-off criterion on
# the "RTL Design" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# "README" sheet (1st tab)
# ---------------------------------------------------------------
$readme = $wb.Worksheets.Item(1)

$readme.Range("B2").Value = 'Copyright 2020, 2023 by the OpenHW Group'
$readme.Range("B3").Value = 'Licensed under the Solderpad Hardware License, Version 2.1'
$readme.Range("B6").Value = 'This document is the Checklist for the TRL-5 Release for OpenHW Group RTL IP.'
$readme.Range("B11").Value = 'What Does “TRL-5” Mean?'
$readme.Range("B12").Value = 'RTL that meets the TRL-5 criteria is complete, functionally correct, validated against a specific software toolchain and ready to be used in a commercial-grade product.  The Specification, Design and Verification are complete and self-consistent.  It has been shown to match the design intent as captured in the specification by means of either dynamic or static verification methods (or both).'
$readme.Range("B14").Value = 'Deliverables: OpenHW provides the following at TRL-5:'
$readme.Range("B22").Value = 'IP that achieves the TRL-5 criteria is subject to a set of RTL design rules and lint checks.  It may or may not have been synthesized and implemented into a physical gate model.'

# ---------------------------------------------------------------
# "RTL Design" sheet (3rd tab)
# ---------------------------------------------------------------
$rtlDesign = $wb.Worksheets.Item(3)

$rtlDesign.Range("D3").Value = 'The version at this release is clearly identified in GitHub and in the release review document.'

# Remove the mailto: hyperlink (and its display text,
# "arjan.bink@silabs.com") that lived in the "Owner" column for this
# row; the cell becomes empty.
$rtlDesign.Range("E3").Hyperlinks.Delete()
$rtlDesign.Range("E3").ClearContents()

$rtlDesign.Range("F9").Value = 'Waiving can be done by applying the WAIVED:<PROJECT_NAME>. Issues labeled with a non-applicable parameter option are waived as well in case the release configuration applies to a different parameter configuration'
